$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "alerttext"
$ws.Range("D2").Value = "Customer added successfully"

$ws.Range("D3").Select()
